$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (everything from the old row 6 "D1, D2" onward
# shifts down by one; the C4-C14/bypass-cap row above stays at row 5).
$ws.Rows(6).Insert()

# Fill in the newly inserted row 6: a 470uF/16V stability capacitor.
$ws.Range("B6").Value = "470uf/16V"
$ws.Range("D6").Value = "additional stability when using pi; place across the 5V and GND pins of SV1"
$ws.Range("C6").Value = 1

# Append two new rows at the bottom of the table for the ribbon cable + IDC
# connectors used to hook the ETA-3400 up to the ET-3400.
$ws.Range("A35").Value = "CAB1"
$ws.Range("D35").Value = "40 pin ribbon cable"
$ws.Range("B35").Value = "40-ribbon"

$ws.Range("A36").Value = "IDC1, IDC2"
$ws.Range("B36").Value = "2x20 IDC conn"
$ws.Range("D36").Value = "2x20 IDC ribbon cable connector"

$ws.Range("F35").Value = "ETA-3400 to ET-3400 cable"
$ws.Range("F36").Value = "ETA-3400 to ET-3400 connectors"

# The IC2 (RAM) row is now row 9 after the insert; update its source part number.
$ws.Range("E9").Value = "digikey 1450-1182-5-ND (use the 70ns version, not the 55ns)"
